$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.17%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.80%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.216"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.89%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07665"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.52%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "'FTXToken"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'1.639"
$ws.Range("D6").Style = "Normal"
$ws.Range("B7").Value = "'MXToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'0.9147"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.17%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'BTSEToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'2.429"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.24%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.1208"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'9.38%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1819"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.05%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.09179"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.61%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.04153"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.93%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.1050"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.21%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.001256"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.32%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'TigerCash"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.005825"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.67%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'UpBots"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.007509"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2,395.62%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.340"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.39%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'GateToken"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'4.312"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.42%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'7.402"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'13.26%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.31%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2714"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.20%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.03996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.49%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'2.65%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004382"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'6.64%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.11%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02484"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'2.78%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05331"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.70%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007847"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.74%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1314"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.71%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006507"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.38%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001912"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.008259"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-5.88%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3346"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.38%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006718"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-3.37%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.12%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.3912"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1,125.06%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003104"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-26.10%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.12%"
$ws.Range("E51").Style = "Normal"
